$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the caption for the "direct_indirect" plot (row 4, column C):
# the year referenced moves from 2018 to 2019, a sentence about Annex B is
# inserted, and a closing sentence about GWP-100 reporting is appended.
$ws.Range("C4").Value = "The stacked bar on the left indicates total global greenhouse gas emissions in 2019, split by sectors based on direct (scope 1) emissions accounting. The arrows shown next to the electricity and heat sector depict the reallocation of these emissions to final sectors as indirect (scope 2) emissions. This increases the contribution to global emissions from the industry and buildings sector (central stacked bar). This reallocation does not imply full lifecycle emissions – see Annex B {A.B.8} for more details. The stacked bar on the far right indicates the shares of subsectors in global emissions when indirect emissions are included. GHG emissions are reported in GtCO2-eq, based on global warming potentials with a 100-year time horizon (GWP-100) from the IPCC Sixth Assessment Working Group 1 Report."

# Move selection back to the top of the sheet and select C5 (mirrors the
# scrolled-back view recorded in the workbook after the edit).
$ws.Application.Goto($ws.Range("A1"))
$ws.Range("C5").Select()
